# Export with no is_pref and no lev distance
#
# The "id" column (B) is regenerated directly from the "speaker_variant"
# column (C) as "#" + lowercase(speaker_variant), instead of being matched
# to some previously-picked / levenshtein-nearest preferred id. The
# "is_prefered" column (D) is no longer populated at all.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> speaker_variant (column C) value for the new export order
$rows = [ordered]@{
    2  = "Octa"
    3  = "Tyter"
    4  = "Octavio"
    5  = "Theoph"
    6  = "Guyd"
    7  = "Brand"
    8  = "Otto"
    9  = "Octav"
    10 = "galdra"
    11 = "Billinc"
    12 = "Keyser"
    13 = "Thes"
    14 = "Galdra"
    15 = "Galdrad"
    16 = "guyd"
    17 = "Thessalia"
    18 = "Echo"
    19 = "Flavio"
    20 = "Camillo"
    21 = "Laura"
    22 = "Thess"
    23 = "billinc"
}

foreach ($r in $rows.Keys) {
    $speakerVariant = $rows[$r]
    $id = "#" + $speakerVariant.ToLower()

    $ws.Cells.Item($r, 2).Value = $id             # id
    $ws.Cells.Item($r, 3).Value = $speakerVariant # speaker_variant
    $ws.Cells.Item($r, 4).Value = ""              # is_prefered (cleared)
}
